$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9349475502967834
$ws.Range("B1").Value = 1.684190273284912
$ws.Range("C1").Value = 4.525660991668701
$ws.Range("D1").Value = 2.200752258300781
$ws.Range("E1").Value = 0.9746428728103638
